$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9569.571
$ws.Range("J32").Value = 8275.25
$ws.Range("L32").Value = 8275.25
$ws.Range("N32").Value = -8927.25

$ws.Range("H38").Value = 221.2
$ws.Range("I38").Value = 221.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 663.5999999999999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -291.5999999999999
$ws.Range("N38").ClearContents()

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").ClearContents()
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

$ws.Range("H112").Value = 1716.4546
$ws.Range("J112").Value = 1717.05
$ws.Range("L112").Value = 5151.15
$ws.Range("N112").Value = -7367.15

$ws.Range("H137").Value = 2908.913
$ws.Range("I137").Value = 1767.0769
$ws.Range("J137").Value = 4393.3
$ws.Range("K137").Value = 5301.2307
$ws.Range("L137").Value = 13179.9
$ws.Range("M137").Value = -2751.2307
$ws.Range("N137").Value = -18279.9


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4565.6665
$ws.Range("I32").Value = 2590.2974
$ws.Range("K32").Value = 2590.2974
$ws.Range("M32").Value = -2303.2974

$ws.Range("H61").Value = 4219.636
$ws.Range("I61").Value = 3774.2856
$ws.Range("K61").Value = 3774.2856
$ws.Range("M61").Value = -3562.2856

$ws.Range("H74").Value = 2672.318
$ws.Range("J74").Value = 2933.6667
$ws.Range("L74").Value = 2933.6667
$ws.Range("N74").Value = -4681.6667

$ws.Range("H77").Value = 2672.318
$ws.Range("J77").Value = 2933.6667
$ws.Range("L77").Value = 14668.3335
$ws.Range("N77").Value = -23404.3335

$ws.Range("H110").Value = 2392.875
$ws.Range("I110").Value = 1857.5
$ws.Range("K110").Value = 1857.5
$ws.Range("M110").Value = 187.5

$ws.Range("H136").Value = 4219.636
$ws.Range("I136").Value = 3774.2856
$ws.Range("K136").Value = 11322.8568
$ws.Range("M136").Value = -8772.856800000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 32557.334
$ws.Range("I74").Value = 39000
$ws.Range("K74").Value = 39000
$ws.Range("M74").Value = -38064

$ws.Range("H77").Value = 32557.334
$ws.Range("I77").Value = 39000
$ws.Range("K77").Value = 117000
$ws.Range("M77").Value = -112320

$ws.Range("H134").Value = 4049.889
$ws.Range("I134").Value = 3558.4285
$ws.Range("K134").Value = 10675.2855
$ws.Range("M134").Value = -8140.2855


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4942.5557
$ws.Range("I105").Value = 3117
$ws.Range("K105").Value = 3117
$ws.Range("M105").Value = -1370


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 253.33333
$ws.Range("J23").Value = 206.88889
$ws.Range("L23").Value = 620.6666700000001
$ws.Range("N23").Value = -1090.66667

$ws.Range("H51").Value = 2172.1428
$ws.Range("I51").Value = 2159.6667
$ws.Range("J51").Value = 2181.5
$ws.Range("K51").Value = 6479.000100000001
$ws.Range("L51").Value = 6544.5
$ws.Range("M51").Value = -6019.000100000001
$ws.Range("N51").Value = -7464.5

$ws.Range("H80").Value = 13238.714
$ws.Range("I80").Value = 18776.111
$ws.Range("J80").Value = 9085.666999999999
$ws.Range("K80").Value = 56328.333
$ws.Range("L80").Value = 27257.001
$ws.Range("M80").Value = -55392.333
$ws.Range("N80").Value = -29129.001

$ws.Range("H83").Value = 13238.714
$ws.Range("I83").Value = 18776.111
$ws.Range("J83").Value = 9085.666999999999
$ws.Range("K83").Value = 168984.999
$ws.Range("L83").Value = 81771.003
$ws.Range("M83").Value = -164304.999
$ws.Range("N83").Value = -91131.003

$ws.Range("H94").Value = 3674.6667
$ws.Range("I94").Value = 3012
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 9036
$ws.Range("L94").Value = 15000
$ws.Range("M94").Value = -8360
$ws.Range("N94").Value = -16352

$ws.Range("H110").Value = 499
$ws.Range("I110").Value = 499
$ws.Range("K110").Value = 1497
$ws.Range("M110").Value = 2593

$ws.Range("H111").Value = 495
$ws.Range("I111").Value = 495
$ws.Range("K111").Value = 1485
$ws.Range("M111").Value = 1582

$ws.Range("H112").Value = 849.5
$ws.Range("I112").Value = 499
$ws.Range("J112").Value = 1200
$ws.Range("K112").Value = 1497
$ws.Range("L112").Value = 3600
$ws.Range("M112").Value = -389
$ws.Range("N112").Value = -5816

$ws.Range("H114").Value = 1597.6666
$ws.Range("I114").Value = 1597.6666
$ws.Range("K114").Value = 4792.9998
$ws.Range("M114").Value = -1538.9998

$ws.Range("H115").Value = 1940.6666
$ws.Range("I115").Value = 1940.6666
$ws.Range("K115").Value = 5821.9998
$ws.Range("M115").Value = -4646.9998

$ws.Range("H117").Value = 839.8
$ws.Range("J117").Value = 499
$ws.Range("L117").Value = 1497
$ws.Range("N117").Value = -8381

$ws.Range("H118").Value = 2229
$ws.Range("I118").Value = 2229
$ws.Range("K118").Value = 6687
$ws.Range("M118").Value = -5444

$ws.Range("H120").Value = 10482.5
$ws.Range("I120").Value = 965
$ws.Range("K120").Value = 2895
$ws.Range("M120").Value = 1943

$ws.Range("H132").Value = 1552.2858
$ws.Range("I132").Value = 1571.7693
$ws.Range("K132").Value = 14145.9237
$ws.Range("M132").Value = -11615.9237

$ws.Range("H141").Value = 1643.3334
$ws.Range("I141").Value = 1643.3334
$ws.Range("K141").Value = 4930.0002
$ws.Range("M141").Value = 249.9997999999996


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 932.4167
$ws.Range("I22").Value = 3099.5
$ws.Range("K22").Value = 3099.5
$ws.Range("M22").Value = -2804.5

$ws.Range("H27").Value = 932.4167
$ws.Range("I27").Value = 3099.5
$ws.Range("K27").Value = 3099.5
$ws.Range("M27").Value = -2992.5

$ws.Range("H46").Value = 2912.4285
$ws.Range("I46").Value = 2297.25
$ws.Range("K46").Value = 2297.25
$ws.Range("M46").Value = -2109.25

$ws.Range("H132").Value = 4550.8423
$ws.Range("I132").Value = 4499
$ws.Range("K132").Value = 13497
$ws.Range("M132").Value = -10967

$ws.Range("H136").Value = 3037.625
$ws.Range("I136").Value = 3059.6
$ws.Range("K136").Value = 9178.799999999999
$ws.Range("M136").Value = -6628.799999999999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 57714.145
$ws.Range("J82").Value = 65000
$ws.Range("L82").Value = 65000
$ws.Range("N82").Value = -65766

$ws.Range("H85").Value = 57714.145
$ws.Range("J85").Value = 65000
$ws.Range("L85").Value = 65000
$ws.Range("N85").Value = -67652

$ws.Range("H107").Value = 477.2857
$ws.Range("I107").Value = 438.6
$ws.Range("J107").Value = 574
$ws.Range("K107").Value = 1315.8
$ws.Range("L107").Value = 1722
$ws.Range("M107").Value = 604.1999999999998
$ws.Range("N107").Value = -5562

